# Weekly update: a new price-report row for "Arveja Verde" at
# "Vega Modelo de Temuco" is inserted as row 24 (pushing the existing
# rows 24-71 down to 25-72). The inserted row carries the new reading;
# everything else keeps its previous values, just shifted one row down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24; this shifts the old rows 24..71 to 25..72
# and copies formatting (incl. the date number format) from the row below.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new data point.
$ws.Range("A24").Value2 = 10
$ws.Range("B24").Value2 = "Vega Modelo de Temuco"
$ws.Range("C24").Value2 = "La Araucanía"
$ws.Range("D24").Value2 = 44536
$ws.Range("E24").Value2 = 9
$ws.Range("F24").Value2 = 100112022
$ws.Range("G24").Value2 = "Arveja Verde"
$ws.Range("H24").Value2 = "Sin especificar"
$ws.Range("I24").Value2 = "Primera"
$ws.Range("J24").Value2 = 290
$ws.Range("K24").Value2 = 13000
$ws.Range("L24").Value2 = 15000
$ws.Range("M24").Value2 = 14138
$ws.Range("N24").Value2 = "$/saco 25 kilos"
$ws.Range("O24").Value2 = "Región de La Araucanía"
$ws.Range("P24").Value2 = 566
$ws.Range("Q24").Value2 = 25
$ws.Range("R24").Value2 = "Hortaliza"
